# Gonzales-09302024-10052024.docx edit: "Add files via upload"
#
# 1. Clear out the "IT Support Intern" entry that was filled in under the
#    "Department Assigned:" table cell (the run is removed entirely, leaving
#    the paragraph empty but keeping its pPr/rPr formatting).
# 2. Remove the now-stale "_GoBack" bookmark left over from the last edit
#    position (Word drops this automatically once the doc is clean-saved).

$d = $word.ActiveDocument

# --- 1. Remove the "IT Support Intern" run -------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "IT Support Intern",   # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                     # Wrap (wdFindContinue)
    $false,                # Format
    "",                    # ReplaceWith
    2                      # Replace (wdReplaceAll)
) | Out-Null

# --- 2. Delete the "_GoBack" bookmark -------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
